$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new daily row (row 63) that the diff adds beneath the
# existing data (last existing data row is 62).
# Column A holds a date-like string ("2025/10/05") that must stay literal
# text (matching rows 2-62), not get auto-converted to a date serial by
# Excel's smart input parsing, hence the leading apostrophe (quote-prefix)
# which forces text entry without altering the number format.
$ws.Range("A63").Value = "'2025/10/05"
$ws.Range("B63").Value = "日"
$ws.Range("C63").Value = 8
$ws.Range("D63").Value = 5
